# Auto-generated edit script: updates currentAveragePrice / Leve price & profit
# columns (H-N) across several Asura_Profits sheets per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1625.591
$ws.Range("I62").Value = 1173.5
$ws.Range("J62").Value = 2002.3334
$ws.Range("K62").Value = 1173.5
$ws.Range("L62").Value = 2002.3334
$ws.Range("M62").Value = -549.5
$ws.Range("N62").Value = -3250.3334

# Row 65
$ws.Range("H65").Value = 1625.591
$ws.Range("I65").Value = 1173.5
$ws.Range("J65").Value = 2002.3334
$ws.Range("K65").Value = 5867.5
$ws.Range("L65").Value = 10011.667
$ws.Range("M65").Value = -2747.5
$ws.Range("N65").Value = -16251.667

# Row 88
$ws.Range("H88").Value = 3358.6667
$ws.Range("I88").Value = 2733.3333
$ws.Range("J88").Value = 3515
$ws.Range("K88").Value = 2733.3333
$ws.Range("L88").Value = 3515
$ws.Range("M88").Value = -2327.3333
$ws.Range("N88").Value = -4327

# Row 91
$ws.Range("H91").Value = 3358.6667
$ws.Range("I91").Value = 2733.3333
$ws.Range("J91").Value = 3515
$ws.Range("K91").Value = 2733.3333
$ws.Range("L91").Value = 3515
$ws.Range("M91").Value = -1329.3333
$ws.Range("N91").Value = -6323

# Row 129
$ws.Range("H129").Value = 1050.1765
$ws.Range("I129").Value = 615.625
$ws.Range("J129").Value = 1183.8846
$ws.Range("K129").Value = 1846.875
$ws.Range("L129").Value = 3551.6538
$ws.Range("M129").Value = 3153.125
$ws.Range("N129").Value = -13551.6538

# Row 135
$ws.Range("H135").Value = 734.81134
$ws.Range("I135").Value = 596.70215
$ws.Range("J135").Value = 1816.6666
$ws.Range("K135").Value = 5370.31935
$ws.Range("L135").Value = 16349.9994
$ws.Range("M135").Value = -2835.31935
$ws.Range("N135").Value = -21419.9994

# Row 137
$ws.Range("H137").Value = 1387.9474
$ws.Range("I137").Value = 1202.7273
$ws.Range("K137").Value = 3608.1819
$ws.Range("M137").Value = -1058.1819

# Row 138
$ws.Range("H138").Value = 4602.593
$ws.Range("I138").Value = 2431.1052
$ws.Range("J138").Value = 9759.875
$ws.Range("K138").Value = 7293.3156
$ws.Range("L138").Value = 29279.625
$ws.Range("M138").Value = -2153.3156
$ws.Range("N138").Value = -39559.625

# Row 141
$ws.Range("H141").Value = 4508.524
$ws.Range("I141").Value = 2150.7368
$ws.Range("J141").Value = 26907.5
$ws.Range("K141").Value = 6452.2104
$ws.Range("L141").Value = 80722.5
$ws.Range("M141").Value = -1272.2104
$ws.Range("N141").Value = -91082.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9153.848
$ws.Range("I32").Value = 10560.525
$ws.Range("J32").Value = 2769.6924
$ws.Range("K32").Value = 10560.525
$ws.Range("L32").Value = 2769.6924
$ws.Range("M32").Value = -10273.525
$ws.Range("N32").Value = -3343.6924

# Row 45
$ws.Range("H45").Value = 1175.9166
$ws.Range("I45").Value = 889
$ws.Range("J45").Value = 1749.75
$ws.Range("K45").Value = 889
$ws.Range("L45").Value = 1749.75
$ws.Range("M45").Value = -512
$ws.Range("N45").Value = -2503.75

# Row 61
$ws.Range("H61").Value = 1481.2812
$ws.Range("I61").Value = 1281.8636
$ws.Range("J61").Value = 1920
$ws.Range("K61").Value = 1281.8636
$ws.Range("L61").Value = 1920
$ws.Range("M61").Value = -1069.8636
$ws.Range("N61").Value = -2344

# Row 74
$ws.Range("H74").Value = 857.8421
$ws.Range("I74").Value = 789.35486
$ws.Range("J74").Value = 1161.1428
$ws.Range("K74").Value = 789.35486
$ws.Range("L74").Value = 1161.1428
$ws.Range("M74").Value = 84.64513999999997
$ws.Range("N74").Value = -2909.1428

# Row 77
$ws.Range("H77").Value = 857.8421
$ws.Range("I77").Value = 789.35486
$ws.Range("J77").Value = 1161.1428
$ws.Range("K77").Value = 3946.7743
$ws.Range("L77").Value = 5805.714
$ws.Range("M77").Value = 421.2257
$ws.Range("N77").Value = -14541.714

# Row 88
$ws.Range("H88").Value = 2337
$ws.Range("I88").Value = 1578.6666
$ws.Range("J88").Value = 2716.1667
$ws.Range("K88").Value = 1578.6666
$ws.Range("L88").Value = 2716.1667
$ws.Range("M88").Value = -1172.6666
$ws.Range("N88").Value = -3528.1667

# Row 91
$ws.Range("H91").Value = 2337
$ws.Range("I91").Value = 1578.6666
$ws.Range("J91").Value = 2716.1667
$ws.Range("K91").Value = 1578.6666
$ws.Range("L91").Value = 2716.1667
$ws.Range("M91").Value = -174.6666
$ws.Range("N91").Value = -5524.1667

# Row 132
$ws.Range("H132").Value = 1991.6875
$ws.Range("I132").Value = 1305.091
$ws.Range("K132").Value = 3915.273
$ws.Range("M132").Value = -1385.273

# Row 135
$ws.Range("H135").Value = 39962.082
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 39958.637
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 39958.637
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -50098.637

# Row 136
$ws.Range("H136").Value = 1481.2812
$ws.Range("I136").Value = 1281.8636
$ws.Range("J136").Value = 1920
$ws.Range("K136").Value = 3845.5908
$ws.Range("L136").Value = 5760
$ws.Range("M136").Value = -1295.5908
$ws.Range("N136").Value = -10860

$ws = $wb.Worksheets.Item("BSM")
# Row 46
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -8096

# Row 123
$ws.Range("H123").Value = 15725
$ws.Range("J123").Value = 15725
$ws.Range("L123").Value = 15725
$ws.Range("N123").Value = -25525

# Row 134
$ws.Range("H134").Value = 2185.1282
$ws.Range("I134").Value = 1848.5807
$ws.Range("J134").Value = 3489.25
$ws.Range("K134").Value = 5545.742099999999
$ws.Range("L134").Value = 10467.75
$ws.Range("M134").Value = -3010.742099999999
$ws.Range("N134").Value = -15537.75

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 6239.0586
$ws.Range("I22").Value = 7929.5386
$ws.Range("J22").Value = 745
$ws.Range("K22").Value = 7929.5386
$ws.Range("L22").Value = 745
$ws.Range("M22").Value = -7579.5386
$ws.Range("N22").Value = -1445

# Row 132
$ws.Range("H132").Value = 357100.5
$ws.Range("I132").Value = 451634
$ws.Range("J132").Value = 2599.875
$ws.Range("K132").Value = 1354902
$ws.Range("L132").Value = 7799.625
$ws.Range("M132").Value = -1352372
$ws.Range("N132").Value = -12859.625

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 12672441
$ws.Range("J131").Value = 14098745
$ws.Range("L131").Value = 42296235
$ws.Range("N131").Value = -42306315

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 24999.908
$ws.Range("I5").Value = 25000
$ws.Range("J5").Value = 24999.9
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 24999.9
$ws.Range("M5").Value = -24888
$ws.Range("N5").Value = -25223.9

# Row 80
$ws.Range("H80").Value = 4000.625
$ws.Range("I80").Value = 4000.7144
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 4000.7144
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -3002.7144
$ws.Range("N80").Value = -5996

# Row 83
$ws.Range("H83").Value = 4000.625
$ws.Range("I83").Value = 4000.7144
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 20003.572
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -15011.572
$ws.Range("N83").Value = -29984

# Row 132
$ws.Range("H132").Value = 1594.0834
$ws.Range("I132").Value = 983.6
$ws.Range("J132").Value = 2611.5557
$ws.Range("K132").Value = 2950.8
$ws.Range("L132").Value = 7834.6671
$ws.Range("M132").Value = -420.8000000000002
$ws.Range("N132").Value = -12894.6671

# Row 134
$ws.Range("H134").Value = 86178.625
$ws.Range("J134").Value = 86178.625
$ws.Range("L134").Value = 258535.875
$ws.Range("N134").Value = -263605.875

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2630.0244
$ws.Range("I136").Value = 2643.9253
$ws.Range("J136").Value = 2567.9333
$ws.Range("K136").Value = 7931.7759
$ws.Range("L136").Value = 7703.7999
$ws.Range("M136").Value = -5381.7759
$ws.Range("N136").Value = -12803.7999

$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 92
$ws.Range("H92").Value = 30320
$ws.Range("J92").Value = 30320
$ws.Range("L92").Value = 30320
$ws.Range("N92").Value = -35312

# Row 94
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
